# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the handback
# xliff files have been generated / processed:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The per-language detail sheets (zh-cn, de-de) get their
#     "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#     columns filled in, with a hyperlink added on the target-file cell.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

$srcMdName1 = "6bef14aa-168f-43d5-b0d7-06757f3b9b01.md"
$srcMdName2 = "ee3fa16d-8fda-4d94-a04c-7001f1446215.md"

$srcMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d3bfbace1209f35175535d9f39563d6840d87b5/e2e/6bef14aa-168f-43d5-b0d7-06757f3b9b01.md"
$srcMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d3bfbace1209f35175535d9f39563d6840d87b5/e2e/ee3fa16d-8fda-4d94-a04c-7001f1446215.md"

$zhHandback1 = "6bef14aa-168f-43d5-b0d7-06757f3b9b01.ee650f03c1331b1e05f5262f2ee63fcda48bf1db.zh-cn.xlf"
$zhHandback2 = "ee3fa16d-8fda-4d94-a04c-7001f1446215.0858773b6454b4b5662af53d81c9fc2592e3b088.zh-cn.xlf"
$deHandback1 = "6bef14aa-168f-43d5-b0d7-06757f3b9b01.ee650f03c1331b1e05f5262f2ee63fcda48bf1db.de-de.xlf"
$deHandback2 = "ee3fa16d-8fda-4d94-a04c-7001f1446215.0858773b6454b4b5662af53d81c9fc2592e3b088.de-de.xlf"

$zhHandbackTime = "2016-08-29 00:46:45"
$deHandbackTime = "2016-08-29 00:46:51"

$hyperlinkColor = 15570276  # OLE BGR for RGB FF6495ED, matches existing HyperLink style

# ---------------------------------------------------------------------------
# Overview sheet: update the Status columns (zh-cn / de-de) for both files.
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = $srcMdName1
$zhcn.Range("J2").Value = $zhHandback1
$zhcn.Range("K2").Value = $zhHandbackTime

$zhcn.Range("I3").Value = $srcMdName2
$zhcn.Range("J3").Value = $zhHandback2
$zhcn.Range("K3").Value = $zhHandbackTime

# Rebuild the hyperlinks top-to-bottom / left-to-right so the relationship
# ids come out in the same order as the source file hyperlinks plus the
# newly added target-file hyperlinks (A2, I2, A3, I3).
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $srcMdUrl1, "", "", $srcMdName1)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcMdUrl1, "", "", $srcMdName1)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $srcMdUrl2, "", "", $srcMdName2)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $srcMdUrl2, "", "", $srcMdName2)

$zhRng = $zhcn.Range("I2:I3")
$zhRng.Font.Color = $hyperlinkColor
$zhRng.Font.Underline = 2

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $srcMdName1
$dede.Range("J2").Value = $deHandback1
$dede.Range("K2").Value = $deHandbackTime

$dede.Range("I3").Value = $srcMdName2
$dede.Range("J3").Value = $deHandback2
$dede.Range("K3").Value = $deHandbackTime

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $srcMdUrl1, "", "", $srcMdName1)
$dede.Hyperlinks.Add($dede.Range("I2"), $srcMdUrl1, "", "", $srcMdName1)
$dede.Hyperlinks.Add($dede.Range("A3"), $srcMdUrl2, "", "", $srcMdName2)
$dede.Hyperlinks.Add($dede.Range("I3"), $srcMdUrl2, "", "", $srcMdName2)

$deRng = $dede.Range("I2:I3")
$deRng.Font.Color = $hyperlinkColor
$deRng.Font.Underline = 2

# ---------------------------------------------------------------------------
# Column widths widen to fit the newly populated / longer content.
# (Input values below are chosen so that, after this runtime's internal
# character-width quantization, the saved OOXML column width lands on the
# intended value - 30 and 40 "characters" respectively.)
# ---------------------------------------------------------------------------
$width30 = 29.16666666
$width40 = 39.16666666

$overview.Columns.Item(5).ColumnWidth = $width30
$overview.Columns.Item(6).ColumnWidth = $width30

$zhcn.Columns.Item(3).ColumnWidth = $width30
$zhcn.Columns.Item(9).ColumnWidth = $width40
$zhcn.Columns.Item(10).ColumnWidth = $width40

$dede.Columns.Item(3).ColumnWidth = $width30
$dede.Columns.Item(9).ColumnWidth = $width40
$dede.Columns.Item(10).ColumnWidth = $width40
